$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E74").Value = (Get-Date -Year 2019 -Month 1 -Day 2).Date
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = "Stunden"
$ws.Range("H74").Value = "Programmieren"
$ws.Range("I74").Value = "Weitere bearbeitung des GUI und erstellung des neuen Fensters zur einzelnen Torrent-Auswahl"

$ws.Range("E75").Value = (Get-Date -Year 2019 -Month 1 -Day 3).Date
$ws.Range("F75").Value = 2
$ws.Range("G75").Value = "Stunden"
$ws.Range("H75").Value = "Programmieren"
$ws.Range("I75").Value = "Implementierung des neuen Fensters und Fehlerbehandlung"

$ws.Range("E76").Value = (Get-Date -Year 2019 -Month 1 -Day 4).Date
$ws.Range("F76").Value = 0.5
$ws.Range("G76").Value = "Stunden"
$ws.Range("H76").Value = "Programmieren"
$ws.Range("I76").Value = "Definition mathematischer Operationen und Code Review"

$ws.Range("H75:H76").Select()
